# Cafe de los maestros Eintritt und Kiosk eingetragen
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of kiosk data appended right below the existing table (row 30)
$ws.Range("A30").Value = 45626
$ws.Range("B30").Value = "Spez 1"
$ws.Range("C30").Value = "Rotwein"
$ws.Range("D30").Value = 7

# Copy number formatting from the row above so the new cells match
# (date format for A, currency format for D) without minting new styles.
$ws.Range("A29").Copy()
$ws.Range("A30").PasteSpecial(-4122)
$ws.Range("D29").Copy()
$ws.Range("D30").PasteSpecial(-4122)

# Expand Table1 to include the new row
$tbl = $ws.ListObjects.Item("Table1")
$tbl.Resize($ws.Range("A1:E30"))

# Match the selection left by the author after entering the new data
$ws.Range("D30").Select()
